{"js": "// Fix #300: a run that carried a trailing <w:br w:type=\"page\"/> together\n// with one or more *extra* duplicate page breaks and the following\n// \"PageN\" text is split apart so that:\n//   - every page break lives alone in its own run (preceded by an empty\n//     <w:t/>), and\n//   - the extra/duplicated page breaks that were inflating the page\n//     count are removed, leaving exactly one break per page marker.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that contains the \"PageN\" / page-break markers\n// (identified by its \"Page1\" text) without disturbing any other\n// paragraph in the body.\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const candidate = paragraphs.items[i];\n  if (candidate.text.indexOf(\"Page1\") === 0) {\n    targetParagraph = candidate;\n    break;\n  }\n}\n\nif (targetParagraph) {\n  // Anchor on the \"Page1\" run so the preceding bookmark start/end\n  // elements (and the paragraph's own attributes) are left untouched;\n  // only the run content from \"Page1\" to the end of the paragraph is\n  // replaced.\n  const page1Results = targetParagraph.search(\"Page1\", { matchWholeWord: false });\n  page1Results.load(\"items\");\n  await context.sync();\n\n  const startRange = page1Results.items[0];\n  const endRange = targetParagraph.getRange(\"End\");\n  const contentRange = startRange.expandTo(endRange);\n\n  const replacementOoxml =\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:t>Page1</w:t></w:r>' +\n    '<w:r><w:t/><w:br w:type=\"page\"/></w:r>' +\n    '<w:r><w:t>Page2</w:t></w:r>' +\n    '<w:r><w:t/><w:br w:type=\"page\"/></w:r>' +\n    '<w:r><w:t>Page3</w:t></w:r>' +\n    '<w:r><w:t/><w:br w:type=\"page\"/></w:r>' +\n    '<w:r><w:t>Page4</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n\n  contentRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix #300: a run that carried a trailing <w:br w:type=\"page\"/> together\n# with one or more *extra* duplicate page breaks and the following\n# \"PageN\" text is split apart so that:\n#   - every page break lives alone in its own run (preceded by an empty\n#     <w:t/>), and\n#   - the extra/duplicated page breaks that were inflating the page\n#     count are removed, leaving exactly one break per page marker.\n$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the \"PageN\" / page-break markers\n# (identified by it starting with \"Page1\") without disturbing any other\n# paragraph in the document.\n$targetPara = $null\nforeach ($para in $d.Paragraphs) {\n  if ($para.Range.Text -like \"Page1*\") {\n    $targetPara = $para\n    break\n  }\n}\n\nif ($targetPara -ne $null) {\n  $paraRange = $targetPara.Range\n  $paraEnd = $paraRange.End\n\n  # Anchor on the \"Page1\" run so the preceding bookmark start/end\n  # elements (and the paragraph's own attributes) are left untouched;\n  # only the run content from \"Page1\" to the end of the paragraph is\n  # replaced.\n  $searchRng = $paraRange.Duplicate\n  $searchRng.Find.Execute(\"Page1\") | Out-Null\n  $startPos = $searchRng.Start\n\n  $targetRng = $d.Range($startPos, $paraEnd)\n\n  $replacementOoxml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>Page1</w:t></w:r><w:r><w:t/><w:br w:type=\"page\"/></w:r><w:r><w:t>Page2</w:t></w:r><w:r><w:t/><w:br w:type=\"page\"/></w:r><w:r><w:t>Page3</w:t></w:r><w:r><w:t/><w:br w:type=\"page\"/></w:r><w:r><w:t>Page4</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n  $targetRng.InsertXML($replacementOoxml)\n}\n"}
